# Generate Report for Handback
# Update the Correspond Handoff / Handback datetimes for the
# "69122473-a34a-4d10-b38d-1fdd26acc32e" entry (row 3) on both the
# zh-cn and de-de sheets.

$wb = $excel.ActiveWorkbook

# zh-cn sheet
$wsZh = $wb.Worksheets.Item("zh-cn")
$wsZh.Range("E3").Value = "2016-03-22 05:16:38"
$wsZh.Range("H3").Value = "2016-03-22 05:17:16"

# de-de sheet
$wsDe = $wb.Worksheets.Item("de-de")
$wsDe.Range("E3").Value = "2016-03-22 05:16:46"
$wsDe.Range("H3").Value = "2016-03-22 05:17:29"
